$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.31   # Current Capital
$summary.Range("B6").Value = 66        # Total Trades
$summary.Range("B9").Value = 31.82     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.31     # Capital
$status.Range("D4").Value = 66         # Trades
$status.Range("F4").Value = 0.31       # P&L %
$status.Range("G4").Value = 31.82      # Win Rate %

# --- New closed trade row (#66 / spreadsheet row 67) ---
# Leading apostrophe keeps the date column as literal text (matches the
# other rows, which store dates as plain inline strings, not date serials).
$newRow = @(66, "'2026-02-17", "15:47:35", "MarketMaking", "DOWN", 0.43, 0.432473, "CLOSED", 0.575, 0, 100.31, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(67, $i + 1).Value = $newRow[$i]
    }
}
